$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 30; existing rows 30-54 shift down to 31-55
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with a new weekly price record
$ws.Cells.Item(30, 1).Value = 9
$ws.Cells.Item(30, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(30, 3).Value = "Metropolitana"
$ws.Cells.Item(30, 4).Value = 44778
$ws.Cells.Item(30, 5).Value = 13
$ws.Cells.Item(30, 6).Value = 100112035
$ws.Cells.Item(30, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(30, 8).Value = "Sin especificar"
$ws.Cells.Item(30, 9).Value = "Primera"
$ws.Cells.Item(30, 10).Value = 20
$ws.Cells.Item(30, 11).Value = 18000
$ws.Cells.Item(30, 12).Value = 18000
$ws.Cells.Item(30, 13).Value = 18000
$ws.Cells.Item(30, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(30, 15).Value = "Hijuelas"
$ws.Cells.Item(30, 16).Value = 1200
$ws.Cells.Item(30, 17).Value = 15
$ws.Cells.Item(30, 18).Value = "Hortaliza"
